# venti menu was added
# Fill in the two previously-empty template rows (7 and 8) with the new
# "Venti build" data, fix a typo in J5, and move the selection to E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing cell J5 ("Воччтанова" -> "Восстанова") ---
$ws.Range("J5").Value = "Сила атаки — 1900+, Восстанова — 200% +, МС — 250-300 ед."

# --- Build names (column B) ---
$ws.Range("B7").Value = "Cборка Венти через МС"
$ws.Range("B8").Value = "Сборка Венти через крит. урон"

# --- Row 7: "Cборка Венти через МС" (Venti build via Elemental Mastery) ---
$ws.Range("D7").Value = "HP || В доп. статах МС, восстановление энергии %"
$ws.Range("E7").Value = "Сила атаки || В доп. статах МС, восстановление энергии %"
$ws.Range("F7").Value = "МС || В доп. статах, восстановление энергии %"
$ws.Range("G7").Value = "МС || В доп. статах, восстановление энергии %"
$ws.Range("H7").Value = "МС || В доп. статах, восстановление энергии %"
$ws.Range("J7").Value = "МС — 600-1000 ед., Восстанова — 180-200% (если в команде нет другого Анемо персонажа)"

# --- Row 8: "Сборка Венти через крит. урон" (Venti build via Crit DMG) ---
$ws.Range("D8").Value = "HP || В доп. статах Крит. урон / Шанс крит. попадания , МС,`nВосстановление энергии %,`nСила атаки %"
$ws.Range("E8").Value = "Сила атаки || В доп. статах Крит. урон / Шанс крит. попадания,`nМС,`nВосстановление энергии %,`nСила атаки %"
$ws.Range("F8").Value = "МС / Сила атаки|| В доп. статах Крит. урон / Шанс крит. попадания,`nМС,`nВосстановление энергии %,`nСила атаки %"
$ws.Range("G8").Value = "Бонус Анемо урона % || В доп. статах Крит. урон / Шанс крит. попадания,`nМС,`nВосстановление энергии %,`nСила атаки %"
$ws.Range("H8").Value = "|| В доп. статах Крит. урон / Шанс крит. попадания,`nМС,`nВосстановление энергии %,`nСила атаки %"
$ws.Range("J8").Value = "Сила атаки — 1500-1800 ед., МС — 200-300 ед., Восстанова — 180-200% (если в команде нет другого Анемо персонажа)"

# --- Weapon suggestions (column I) ---
$ws.Range("I7").Value = "Элегия погибели,  Бесструнный, Гаснущие сумерки, Церемониальный лук, Боевой лук Фавония, Охотник во тьме, Ода анемонии "
$ws.Range("I8").Value = "Элегия погибели, Аква симулякрум, Небесное крыло, Бесструнный, Гаснущие сумерки, Церемониальный лук, Боевой лук Фавония, Охотник во тьме, Ода анемонии"

# --- Artifact set suggestions (column C), same text for both rows ---
$ws.Range("C7").Value = "Изумрудная тень, Позолоченные сны, Церемония древней знати, Инструктор, Изгнанник"
$ws.Range("C8").Value = "Изумрудная тень, Позолоченные сны, Церемония древней знати, Инструктор, Изгнанник"

# --- Move the active selection from B7 to E7 ---
$ws.Range("E7").Select() | Out-Null
